$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2023" data column (T), built from the existing "2022"
#     column (S): duplicate its formatting, then its values, and finally
#     overwrite the handful of cells whose 2023 figures differ. ---

$ws.Range("S3:S13").Copy()
$ws.Range("T3:T13").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("S3:S13").Copy()
$ws.Range("T3:T13").PasteSpecial(-4163)   # xlPasteValues

$ws.Cells.Item(3, 20).Value = 2023
$ws.Cells.Item(4, 20).Value = 0
$ws.Cells.Item(5, 20).Value = "-"
$ws.Cells.Item(6, 20).Value = "-"
$ws.Cells.Item(7, 20).Value = "-"
$ws.Cells.Item(8, 20).Value = "-"
$ws.Cells.Item(9, 20).Value = "-"
$ws.Cells.Item(10, 20).Value = "-"
$ws.Cells.Item(11, 20).Value = "-"
$ws.Cells.Item(12, 20).Value = 0.001731197036190674
$ws.Cells.Item(13, 20).Value = "-"

# --- Columns A-C got a little narrower once the new column was added ---
$ws.Columns.Item(1).ColumnWidth = 32.67
$ws.Columns.Item(2).ColumnWidth = 32.67
$ws.Columns.Item(3).ColumnWidth = 32.67

# --- Leave the cursor on A1 instead of the old "T3" selection ---
$ws.Range("A1").Select()
